# Insert a new data row just before the existing "2026/12/29" block.
# Everything from row 696 onward shifts down by one (737 -> 738 rows of data
# overall; sheet dimension grows from A1:D737 to A1:D738).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A696").EntireRow.Insert()

# Column A holds dates stored as plain text (e.g. "2026/01/23"), not real
# date serials, so force the cell to Text format before assigning the
# string -- otherwise Excel auto-converts the "yyyy/mm/dd"-looking text
# into a date serial number. Clear the format afterwards so the new row
# ends up with the same (default/no explicit style) formatting as all the
# other plain data rows.
$newRow = 696
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2026/01/23"
$ws.Cells.Item($newRow, 1).ClearFormats()
$ws.Cells.Item($newRow, 2).Value = "金"
$ws.Cells.Item($newRow, 3).Value = 3
$ws.Cells.Item($newRow, 4).Value = 201
